$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency price (col D) and 1h volume-change (col E) feed.
# A leading apostrophe forces plain-number-looking price strings (e.g.
# "1.00", "70.62") to stay as text instead of being auto-parsed into
# numbers by Excel, matching the source feed's original text formatting.

$ws.Range("D2").Value = "42.215.40"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "2.315.33"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'312.79"
$ws.Range("E5").Value = "  -5.40%  "
$ws.Range("D6").Value = "'105.99"
$ws.Range("E6").Value = "  +5.94%  "
$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  -1.56%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.610"
$ws.Range("E9").Value = "  -2.82%  "
$ws.Range("D10").Value = "'40.19"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").Value = "'8.30"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").Value = "'15.56"
$ws.Range("E15").Value = "  -4.61%  "
$ws.Range("D16").Value = "2.659.66"
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("D17").Value = "2.339.23"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "42.153.11"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").Value = "'7.73"
$ws.Range("E19").Value = "  -2.33%  "
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").Value = "'74.65"
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("E22").Value = "  -6.83%  "
$ws.Range("D23").Value = "'259.25"
$ws.Range("E23").Value = "  -3.53%  "
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").Value = "'9.27"
$ws.Range("E25").Value = "  -7.12%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("D27").Value = "'11.04"
$ws.Range("E27").Value = "  -3.70%  "
$ws.Range("E28").Value = "  +3.42%  "
$ws.Range("D29").Value = "'22.84"
$ws.Range("E29").Value = "  -1.43%  "
$ws.Range("D30").Value = "'35.55"
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("D31").Value = "'0.0894"
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("D32").Value = "'162.95"
$ws.Range("E32").Value = "  -7.34%  "
$ws.Range("E33").Value = "  -5.33%  "
$ws.Range("E34").Value = "  -3.25%  "
$ws.Range("E35").Value = "  -2.57%  "
$ws.Range("E36").Value = "  +11.62%  "
$ws.Range("D37").Value = "'4.51"
$ws.Range("E37").Value = "  -1.57%  "
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").Value = "'2.76"
$ws.Range("E39").Value = "  -7.21%  "
$ws.Range("E40").Value = "  -4.47%  "
$ws.Range("D41").Value = "'98.06"
$ws.Range("E41").Value = "  +7.31%  "
$ws.Range("E42").Value = "  -3.27%  "
$ws.Range("D43").Value = "'70.62"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("E46").Value = "  +2.50%  "
$ws.Range("D47").Value = "'111.78"
$ws.Range("E47").Value = "  -5.01%  "
$ws.Range("D48").Value = "'5.39"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").Value = "'8.94"
$ws.Range("E49").Value = "  -1.98%  "
$ws.Range("D50").Value = "'75.00"
$ws.Range("E50").Value = "  +7.20%  "
$ws.Range("E51").Value = "  +0.00%  "
